# Add a new "2022-Q1" holdings sheet (positioned between "2021-Q4" and
# "总计") and record its summary row in the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet. Clone it from "总计" (placed right
#    after "2021-Q4") so it starts out with the same cell formatting
#    used across the workbook, then rebuild it into the 8-column
#    holdings layout used by the other quarterly sheets.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet.Copy($null, $q4Sheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"

# Extend the header row (currently B:D) out to H, copying D1's format
# onto the new E:H header cells.
$newSheet.Range("D1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fill in the single holder row for 2022-Q1 (A2 already holds 0 with
# the index style carried over from the cloned sheet).
$newSheet.Range("B2").Value = "'005457"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "景顺长城量化小盘股票"
$newSheet.Range("D2").Value = "'9.49"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'93.39"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'1.95"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.1851"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 5

# ------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing 2021-Q4 summary row
#    down to row 3 and add the new 2022-Q1 summary row at row 2.
#    (Re-fetch the sheet by name: the Copy() above re-bound the old
#    $totalSheet variable to the newly created sheet instead of the
#    original "总计" sheet.)
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3"))
$totalSheet.Range("A3").Value = 1

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.19
